$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.319.61'
$ws.Range('E2').Value = '  +2.11%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.658.69'
$ws.Range('E3').Value = '  +1.12%  '
$ws.Range('E4').Value = '  -0.54%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '219.68'
$ws.Range('E5').Value = '  +0.94%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.505'
$ws.Range('E6').Value = '  +0.36%  '
$ws.Range('E7').Value = '  -0.53%  '
$ws.Range('E8').Value = '  +2.02%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.0625'
$ws.Range('E9').Value = '  +0.24%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.94'
$ws.Range('E10').Value = '  +4.54%  '
$ws.Range('E11').Value = '  +0.58%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.890.49'
$ws.Range('E12').Value = '  +1.15%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.653.82'
$ws.Range('E13').Value = '  +1.11%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.19'
$ws.Range('E14').Value = '  +1.07%  '
$ws.Range('E15').Value = '  +0.96%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '67.15'
$ws.Range('E16').Value = '  +3.95%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '27.299.97'
$ws.Range('E17').Value = '  +2.13%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.0₃0733'
$ws.Range('E18').Value = '  +0.64%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '222.06'
$ws.Range('E19').Value = '  +4.35%  '
$ws.Range('E20').Value = '  -0.51%  '
$ws.Range('E21').Value = '  +1.92%  '
$ws.Range('E22').Value = '  +8.18%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.44'
$ws.Range('E23').Value = '  +5.22%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.26'
$ws.Range('E24').Value = '  +0.16%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '147.04'
$ws.Range('E25').Value = '  +1.08%  '
$ws.Range('E26').Value = '  -0.45%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.41'
$ws.Range('E27').Value = '  +3.93%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.120'
$ws.Range('E28').Value = '  +1.50%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '16.03'
$ws.Range('E29').Value = '  +2.55%  '
$ws.Range('E30').Value = '  +1.23%  '
$ws.Range('E31').Value = '  +0.50%  '
$ws.Range('E32').Value = '  +0.08%  '
$ws.Range('E33').Value = '  +0.58%  '
$ws.Range('E34').Value = '  +2.31%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.263.03'
$ws.Range('E35').Value = '  -1.25%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.46'
$ws.Range('E36').Value = '  +0.85%  '
$ws.Range('E37').Value = '  +1.80%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.537'
$ws.Range('E38').Value = '  +0.71%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.838'
$ws.Range('E39').Value = '  +2.90%  '
$ws.Range('E40').Value = '  -0.48%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.819'
$ws.Range('E41').Value = '  +1.82%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.36'
$ws.Range('E42').Value = '  +1.94%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.800.96'
$ws.Range('E43').Value = '  +1.30%  '
$ws.Range('E44').Value = '  -3.83%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '61.71'
$ws.Range('E45').Value = '  +1.25%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '91.93'
$ws.Range('E46').Value = '  +0.70%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.60'
$ws.Range('E47').Value = '  +1.15%  '
$ws.Range('E48').Value = '  -0.86%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.67'
$ws.Range('E49').Value = '  +0.88%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0978'
$ws.Range('E50').Value = '  +1.58%  '
$ws.Range('E51').Value = '  +0.25%  '
